$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.074.67"
$ws.Range("E2").Value = "  +1.88%  "
$ws.Range("D3").Value = "3.142.40"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.19"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.40"
$ws.Range("E6").Value = "  +2.06%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "3.130.65"
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("E10").Value = "  +11.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.74"
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.470"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000254"
$ws.Range("E13").Value = "  +4.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.52"
$ws.Range("E14").Value = "  +5.09%  "
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("D16").Value = "3.662.99"
$ws.Range("E16").Value = "  +1.03%  "
$ws.Range("D17").Value = "63.958.70"
$ws.Range("E17").Value = "  +1.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.18"
$ws.Range("E18").Value = "  -2.40%  "
$ws.Range("D19").Value = "3.141.66"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "468.25"
$ws.Range("E20").Value = "  +3.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.37"
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.737"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.57"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.30"
$ws.Range("E24").Value = "  -3.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.59"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.00"
$ws.Range("E27").Value = "  +8.42%  "
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.87"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.20"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("E33").Value = "  -5.04%  "
$ws.Range("D34").Value = "0.0₃0879"
$ws.Range("E34").Value = "  +9.52%  "
$ws.Range("E35").Value = "  +7.75%  "
$ws.Range("E36").Value = "  +1.23%  "
$ws.Range("E37").Value = "  +13.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.14"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.96"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "456.08"
$ws.Range("E40").Value = "  +6.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.73"
$ws.Range("E41").Value = "  -1.67%  "
$ws.Range("E42").Value = "  +0.58%  "
$ws.Range("D43").Value = "2.905.52"
$ws.Range("E43").Value = "  -1.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.279"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "126.94"
$ws.Range("E47").Value = "  +2.10%  "
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.69"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.80"
$ws.Range("E51").Value = "  +0.19%  "
